# Scheduled market-data refresh: update Leve profit-tracking sheets
# (currentAveragePrice / LevePrice / LeveProfit columns H:N) with the
# latest pulled values, per sheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2400.1428
$ws.Range("I2").Value = 2700
$ws.Range("K2").Value = 2700
$ws.Range("M2").Value = -2587
$ws.Range("H4").Value = 563.5714
$ws.Range("I4").Value = 112.5
$ws.Range("J4").Value = 1165
$ws.Range("K4").Value = 112.5
$ws.Range("L4").Value = 1165
$ws.Range("M4").Value = 1.5
$ws.Range("N4").Value = -1393
$ws.Range("H76").Value = 5976.6665
$ws.Range("I76").Value = 5976.6665
$ws.Range("K76").Value = 5976.6665
$ws.Range("M76").Value = -5661.6665
$ws.Range("H79").Value = 5976.6665
$ws.Range("I79").Value = 5976.6665
$ws.Range("K79").Value = 5976.6665
$ws.Range("M79").Value = -4884.6665
$ws.Range("H137").Value = 92092.05
$ws.Range("I137").Value = 150697.92
$ws.Range("J137").Value = 4183.25
$ws.Range("K137").Value = 452093.76
$ws.Range("L137").Value = 12549.75
$ws.Range("M137").Value = -449543.76
$ws.Range("N137").Value = -17649.75
$ws.Range("H138").Value = 2726.6216
$ws.Range("I138").Value = 1797.7142
$ws.Range("K138").Value = 5393.142599999999
$ws.Range("M138").Value = -253.1425999999992

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1455.7273
$ws.Range("I2").Value = 1268.75
$ws.Range("K2").Value = 1268.75
$ws.Range("M2").Value = -1155.75
$ws.Range("H61").Value = 3545.7273
$ws.Range("I61").Value = 2998.8
$ws.Range("K61").Value = 2998.8
$ws.Range("M61").Value = -2786.8
$ws.Range("H110").Value = 9918.85
$ws.Range("I110").Value = 3266.625
$ws.Range("K110").Value = 3266.625
$ws.Range("M110").Value = -1221.625
$ws.Range("H116").Value = 1455.7273
$ws.Range("I116").Value = 1268.75
$ws.Range("K116").Value = 1268.75
$ws.Range("M116").Value = 1025.25
$ws.Range("H132").Value = 1841.25
$ws.Range("I132").Value = 1187.1482
$ws.Range("K132").Value = 3561.4446
$ws.Range("M132").Value = -1031.4446
$ws.Range("H136").Value = 3545.7273
$ws.Range("I136").Value = 2998.8
$ws.Range("K136").Value = 8996.400000000001
$ws.Range("M136").Value = -6446.400000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1455.7273
$ws.Range("I3").Value = 1268.75
$ws.Range("K3").Value = 1268.75
$ws.Range("M3").Value = -1154.75
$ws.Range("H20").Value = 1916.1333
$ws.Range("I20").Value = 1293.3
$ws.Range("K20").Value = 1293.3
$ws.Range("M20").Value = -1046.3
$ws.Range("H99").Value = 3644.3635
$ws.Range("I99").Value = 2932.5
$ws.Range("K99").Value = 2932.5
$ws.Range("M99").Value = -1434.5
$ws.Range("H134").Value = 6078.263
$ws.Range("I134").Value = 3332.8333
$ws.Range("K134").Value = 9998.499899999999
$ws.Range("M134").Value = -7463.499899999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 489.66666
$ws.Range("I15").Value = 489.66666
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 489.66666
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -319.66666
$ws.Range("N15").ClearContents()
$ws.Range("H57").Value = 20048.8
$ws.Range("I57").Value = 19900
$ws.Range("J57").Value = 20148
$ws.Range("K57").Value = 19900
$ws.Range("L57").Value = 20148
$ws.Range("M57").Value = -19340
$ws.Range("N57").Value = -21268
$ws.Range("H58").Value = 1968.8695
$ws.Range("I58").Value = 1643.3125
$ws.Range("K58").Value = 1643.3125
$ws.Range("M58").Value = -1440.3125
$ws.Range("H136").Value = 1968.8695
$ws.Range("I136").Value = 1643.3125
$ws.Range("K136").Value = 4929.9375
$ws.Range("M136").Value = -2379.9375
$ws.Range("H138").Value = 113499.9
$ws.Range("J138").Value = 186999.8
$ws.Range("L138").Value = 186999.8
$ws.Range("N138").Value = -197279.8
$ws.Range("H140").Value = 89000
$ws.Range("J140").Value = 89000
$ws.Range("L140").Value = 89000
$ws.Range("N140").Value = -99360

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 90
$ws.Range("J20").Value = 90
$ws.Range("L20").Value = 270
$ws.Range("N20").Value = -724
$ws.Range("H23").Value = 270.42856
$ws.Range("I23").Value = 37
$ws.Range("K23").Value = 111
$ws.Range("M23").Value = 124
$ws.Range("H34").Value = 271
$ws.Range("J34").Value = 339.8
$ws.Range("L34").Value = 1019.4
$ws.Range("N34").Value = -1187.4
$ws.Range("H55").Value = 94656.73
$ws.Range("I55").Value = 749.6667
$ws.Range("J55").Value = 129871.875
$ws.Range("K55").Value = 2249.0001
$ws.Range("L55").Value = 389615.625
$ws.Range("M55").Value = -2072.0001
$ws.Range("N55").Value = -389969.625
$ws.Range("H75").Value = 3388.3333
$ws.Range("H78").Value = 3388.3333
$ws.Range("H80").Value = 2068.111
$ws.Range("I80").Value = 2252.75
$ws.Range("J80").Value = 1920.4
$ws.Range("K80").Value = 6758.25
$ws.Range("L80").Value = 5761.200000000001
$ws.Range("M80").Value = -5822.25
$ws.Range("N80").Value = -7633.200000000001
$ws.Range("H83").Value = 2068.111
$ws.Range("I83").Value = 2252.75
$ws.Range("J83").Value = 1920.4
$ws.Range("K83").Value = 20274.75
$ws.Range("L83").Value = 17283.6
$ws.Range("M83").Value = -15594.75
$ws.Range("N83").Value = -26643.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 16671413
$ws.Range("I19").Value = 7500
$ws.Range("J19").Value = 33335326
$ws.Range("K19").Value = 7500
$ws.Range("L19").Value = 33335326
$ws.Range("M19").Value = -7212
$ws.Range("N19").Value = -33335902

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 50429
$ws.Range("J137").Value = 50429
$ws.Range("L137").Value = 50429
$ws.Range("N137").Value = -60629
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 11464.143
$ws.Range("I33").Value = 1500
$ws.Range("J33").Value = 18937.25
$ws.Range("K33").Value = 1500
$ws.Range("L33").Value = 18937.25
$ws.Range("M33").Value = -1250
$ws.Range("N33").Value = -19437.25
$ws.Range("H36").Value = 11464.143
$ws.Range("I36").Value = 1500
$ws.Range("J36").Value = 18937.25
$ws.Range("K36").Value = 1500
$ws.Range("L36").Value = 18937.25
$ws.Range("M36").Value = -1250
$ws.Range("N36").Value = -19437.25
$ws.Range("H81").Value = 1549.75
$ws.Range("I81").Value = 1480
$ws.Range("K81").Value = 2960
$ws.Range("M81").Value = -1899
$ws.Range("H84").Value = 1549.75
$ws.Range("I84").Value = 1480
$ws.Range("K84").Value = 14800
$ws.Range("M84").Value = -9496
$ws.Range("H141").Value = 150000
$ws.Range("J141").Value = 150000
$ws.Range("L141").Value = 150000
